$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily rows to append after the existing last row (row 207)
$newRows = @(
    @{ Row = 208; Serie = "28-10-2021"; B = 2.77; C = 3.68; D = 4.81; E = 5.65; F = 0.16 },
    @{ Row = 209; Serie = "29-10-2021"; B = 2.79; C = 3.75; D = 4.8;  E = 5.69; F = 0.11 },
    @{ Row = 210; Serie = "02-11-2021"; B = 2.79; C = 3.72; D = 4.95; E = 5.75; F = 0.09 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $cellA = $ws.Cells.Item($row, 1)

    if ($r.Serie -eq "02-11-2021") {
        # "02-11-2021" is ambiguous (day <= 12) and Excel's smart-parsing
        # would otherwise silently turn it into a date serial. Force the
        # cell to Text first so it is stored as the literal string, then
        # restore the default "Normal" style so no stray number format is
        # left behind on the cell.
        $cellA.NumberFormat = "@"
        $cellA.Value = $r.Serie
        $cellA.Style = "Normal"
    } else {
        $cellA.Value = $r.Serie
    }

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}
